$wb = $excel.ActiveWorkbook

# --- "hvdc" sheet: replace single "marginal_cost" column with three cost
#     coefficient columns (costc2, costc1, costc0), mirroring the
#     generator sheet's cost model. ---
$hvdc = $wb.Sheets.Item("hvdc")

$hvdc.Range("O1").Value = "costc2"
$hvdc.Range("P1").Value = "costc1"
$hvdc.Range("Q1").Value = "costc0"

$hvdc.Range("P2").Value = 0
$hvdc.Range("Q2").Value = 0

# --- Sheet selection / active tab bookkeeping: the workbook now opens on
#     the "hvdc" sheet (first tab) with Q2 selected, instead of the
#     "timeseries" sheet with H3 selected. ---
$timeseries = $wb.Sheets.Item("timeseries")
$wb.ActiveSheet.Range("W1").Select()

$hvdc.Activate()
$wb.ActiveSheet.Range("Q2").Select()
